$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metrics")

# Update header row (row 1) metric column names to match new evaluation library's
# naming convention (ragas-style metric names) for LightRAG, GraphRAG, NaiveRAG.

$ws.Range("C1").Value = "LightRAG_answer_relevancy"
$ws.Range("D1").Value = "LightRAG_factual_correctness(mode=f1)"
$ws.Range("E1").Value = "LightRAG_factual_correctness(mode=recall)"
$ws.Range("F1").Value = "LightRAG_semantic_similarity"

$ws.Range("G1").Value = "GraphRAG_answer_relevancy"
$ws.Range("H1").Value = "GraphRAG_factual_correctness(mode=f1)"
$ws.Range("I1").Value = "GraphRAG_factual_correctness(mode=recall)"
$ws.Range("J1").Value = "GraphRAG_semantic_similarity"

$ws.Range("K1").Value = "NaiveRAG_answer_relevancy"
$ws.Range("L1").Value = "NaiveRAG_factual_correctness(mode=f1)"
$ws.Range("M1").Value = "NaiveRAG_factual_correctness(mode=recall)"
$ws.Range("N1").Value = "NaiveRAG_semantic_similarity"
